$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "cryptos" price/volume table with the latest values pulled from
# coinranking.com. Price strings that look like plain numbers are written
# with a leading apostrophe so Excel keeps them as text (preserving things
# like trailing zeros, e.g. "0.9990") instead of silently re-parsing them
# into a Double and losing formatting.

$ws.Range("D2").Value = '29.895.68'
$ws.Range("E2").Value = '  -0.54%  '
$ws.Range("D3").Value = '1.900.11'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("D4").Value = '''0.9990'
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").Value = '''0.7628'
$ws.Range("E5").Value = '  +4.69%  '
$ws.Range("D6").Value = '''240.26'
$ws.Range("E6").Value = '  -0.95%  '
$ws.Range("D7").Value = '''0.9991'
$ws.Range("E7").Value = '  -0.23%  '
$ws.Range("D8").Value = '1.896.87'
$ws.Range("E8").Value = '  +0.41%  '
$ws.Range("D9").Value = '''0.3066'
$ws.Range("E9").Value = '  -0.97%  '
$ws.Range("D10").Value = '''25.63'
$ws.Range("E10").Value = '  -1.96%  '
$ws.Range("D11").Value = '''0.06850'
$ws.Range("E11").Value = '  -0.64%  '
$ws.Range("D12").Value = '''0.07958'
$ws.Range("E12").Value = '  +0.19%  '
$ws.Range("D13").Value = '''0.7459'
$ws.Range("E13").Value = '  -3.09%  '
$ws.Range("D14").Value = '1.890.51'
$ws.Range("E14").Value = '  -0.38%  '
$ws.Range("D15").Value = '''5.182'
$ws.Range("E15").Value = '  -1.32%  '
$ws.Range("D16").Value = '''90.96'
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").Value = '29.901.84'
$ws.Range("E17").Value = '  -0.47%  '
$ws.Range("D18").Value = '''13.97'
$ws.Range("E18").Value = '  -1.21%  '
$ws.Range("D19").Value = '''5.968'
$ws.Range("E19").Value = '  +3.88%  '
$ws.Range("D20").Value = '''242.62'
$ws.Range("E20").Value = '  +2.42%  '
$ws.Range("D21").Value = '''0.000007701'
$ws.Range("E21").Value = '  -0.73%  '
$ws.Range("D22").Value = '''0.9990'
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("D23").Value = '''0.9986'
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("D24").Value = '''6.953'
$ws.Range("E24").Value = '  +1.07%  '
$ws.Range("D25").Value = '''166.49'
$ws.Range("E25").Value = '  +0.47%  '
$ws.Range("D26").Value = '''9.241'
$ws.Range("E26").Value = '  -0.74%  '
$ws.Range("D27").Value = '''18.72'
$ws.Range("E27").Value = '  -0.90%  '
$ws.Range("D28").Value = '''0.1288'
$ws.Range("E28").Value = '  +1.17%  '
$ws.Range("D29").Value = '''2.042'
$ws.Range("E29").Value = '  +1.17%  '
$ws.Range("D30").Value = '''1.390'
$ws.Range("E30").Value = '  +2.83%  '
$ws.Range("D31").Value = '''1.514'
$ws.Range("E31").Value = '  -1.43%  '
$ws.Range("D32").Value = '''4.267'
$ws.Range("E32").Value = '  -0.63%  '
$ws.Range("D33").Value = '''4.061'
$ws.Range("E33").Value = '  -0.21%  '
$ws.Range("D34").Value = '''0.05220'
$ws.Range("E34").Value = '  +2.72%  '
$ws.Range("D35").Value = '''1.257'
$ws.Range("E35").Value = '  -0.99%  '
$ws.Range("D36").Value = '''0.7275'
$ws.Range("E36").Value = '  -0.98%  '
$ws.Range("D37").Value = '''2.711'
$ws.Range("E37").Value = '  -0.74%  '
$ws.Range("D38").Value = '''0.01924'
$ws.Range("E38").Value = '  +0.28%  '
$ws.Range("D39").Value = '''2.778'
$ws.Range("E39").Value = '  +0.28%  '
$ws.Range("D40").Value = '''6.156'
$ws.Range("E40").Value = '  -2.83%  '
$ws.Range("D41").Value = '''0.4419'
$ws.Range("E41").Value = '  -0.23%  '
$ws.Range("D42").Value = '''71.99'
$ws.Range("E42").Value = '  -3.53%  '
$ws.Range("D43").Value = '''0.9993'
$ws.Range("E43").Value = '  -0.12%  '
$ws.Range("D44").Value = '''1.888'
$ws.Range("E44").Value = '  -1.90%  '
$ws.Range("D45").Value = '''0.8286'
$ws.Range("E45").Value = '  -0.92%  '
$ws.Range("D46").Value = '''7.642'
$ws.Range("E46").Value = '  +1.17%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '''9.821'
$ws.Range("E47").Value = '  +0.98%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '''99.96'
$ws.Range("E48").Value = '  -0.83%  '
$ws.Range("D49").Value = '2.042.20'
$ws.Range("E49").Value = '  -0.70%  '
$ws.Range("D50").Value = '''36.07'
$ws.Range("E50").Value = '  -3.97%  '
$ws.Range("D51").Value = '''0.05939'
$ws.Range("E51").Value = '  -0.33%  '
